$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (CAS id / group, Corr with BMI, Corr with Age, Corr with HIV, Likely Influenced)
$data = @(
    @("000074-95-3_group1", 0.134, -0.1, 0.048, "Unlikely"),
    @("000098-86-2_group2", 0.08500000000000001, 0.065, 0.007, "Unlikely"),
    @("001640-89-7_group2", -0.053, 0.02, -0.08799999999999999, "Unlikely"),
    @("002216-33-3_group1", -0.038, 0.017, -0.07099999999999999, "Unlikely"),
    @("000098-83-9_group2", -0.067, 0.08400000000000001, -0.017, "Unlikely"),
    @("000124-19-6_group1", 0.051, 0.033, -0.018, "Unlikely"),
    @("000620-14-4_group1", 0.08699999999999999, 0.092, -0.045, "Unlikely"),
    @("000095-16-9_group2", -0.106, 0.031, -0.044, "Unlikely"),
    @("000115-11-7_group2", -0.019, 0.099, -0.058, "Unlikely"),
    @("000617-94-7_group1", -0.114, 0.008, 0.039, "Unlikely")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row++
}
